# ShopSpellTable: remove the "type|String" column (column B).
# Before: id|String, type|String, count|Int, price|Int
# After : id|String, count|Int, price|Int
$wb = $excel.ActiveWorkbook

$shopSpell = $wb.Worksheets.Item("ShopSpellTable")
$shopSpell.Columns.Item(2).Delete()

# ShopActorTable: remove the "type|String" column (column B) and drop the
# trailing "Actor50" row (the table now only lists the 1x and 10x pulls).
# Before: id|String, type|String, count|Int, price|Int  (4 rows of data)
# After : id|String, count|Int, price|Int                (2 rows of data)
$shopActor = $wb.Worksheets.Item("ShopActorTable")
$shopActor.Columns.Item(2).Delete()
$shopActor.Rows.Item(4).Delete()

# ShopActorTable becomes the active/selected sheet (previously it was
# GachaActorTable).
$shopActor.Activate()
